{"js": "// Remove the trailing \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" copyright\n// line right after it, and the blank paragraph that precedes them \u2014 i.e.\n// the footer block that used to sit right after the \"Requisitos\" list.\n\nconst body = context.document.body;\n\n// Locate the \"Ver no Jupiter...\" paragraph via search so we don't depend on\n// hard-coded paragraph indices.\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const hit = results.items[0];\n  const paras = hit.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n\n  const jupiterPara = paras.items[0];\n  const copyrightPara = jupiterPara.getNext();\n  const blankPara = jupiterPara.getPrevious();\n\n  // Delete from bottom to top so earlier references stay valid.\n  copyrightPara.delete();\n  jupiterPara.delete();\n  blankPara.delete();\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" copyright\n# line right after it, and the blank paragraph that precedes them -- i.e.\n# the footer block that used to sit right after the \"Requisitos\" list.\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\n\nif ($found) {\n    $hitStart = $searchRange.Start\n\n    # Work out which 1-based Paragraphs index contains the hit by comparing\n    # offsets (Paragraphs.Item(...).Range.Text is unreliable on ranges coming\n    # from enumeration, but Start/End are trustworthy).\n    $count = $d.Paragraphs.Count\n    $targetIndex = -1\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $r = $p.Range\n        if ($hitStart -ge $r.Start -and $hitStart -lt $r.End) {\n            $targetIndex = $i\n        }\n    }\n\n    if ($targetIndex -gt 0) {\n        # Delete the next paragraph (copyright notice), then the matched\n        # paragraph itself, then the preceding blank paragraph -- bottom to\n        # top so earlier indices stay valid while later ones are removed.\n        $d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n        $d.Paragraphs.Item($targetIndex).Range.Delete()\n        $d.Paragraphs.Item($targetIndex - 1).Range.Delete()\n    }\n}\n"}
